$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2364713333333333
$ws.Range("H2").Value = 0.709414
$ws.Range("I2").Value = 0.002249544876489787
$ws.Range("J2").Value = 0.002249544876489787
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.1013854288673333
$ws.Range("R2").Value = 0.912468859806
$ws.Range("S2").Value = 0.0000092692585462003356797
$ws.Range("T2").Value = 0.0000092692585462003356797
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2364713333333333
$ws.Range("H3").Value = 0.709414
$ws.Range("I3").Value = 0.002249544876489787
$ws.Range("J3").Value = 0.002249544876489787
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 18.971205231706
$ws.Range("R3").Value = 170.740847085354
$ws.Range("S3").Value = 0.00173446034790479
$ws.Range("T3").Value = 0.00173446034790479
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2364713333333333
$ws.Range("H4").Value = 0.709414
$ws.Range("I4").Value = 0.002249544876489787
$ws.Range("J4").Value = 0.002249544876489787
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 5.532513504173556
$ws.Range("R4").Value = 49.792621537562
$ws.Range("S4").Value = 0.0005058152700387974
$ws.Range("T4").Value = 0.0005058152700387974
$ws.Range("I5").Value = 0.9862688099613843
$ws.Range("J5").Value = 0.9862688099613843
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 44.450451876488
$ws.Range("R5").Value = 400.054066888392
$ws.Range("S5").Value = 0.004063924525858152
$ws.Range("T5").Value = 0.004063924525858152
$ws.Range("I6").Value = 0.9862688099613843
$ws.Range("J6").Value = 0.9862688099613843
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 8317.552676079191
$ws.Range("R6").Value = 74857.97408471271
$ws.Range("S6").Value = 0.7604401055215099
$ws.Range("T6").Value = 0.7604401055215099
$ws.Range("I7").Value = 0.9862688099613843
$ws.Range("J7").Value = 0.9862688099613843
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 2425.62198553291
$ws.Range("R7").Value = 21830.59786979618
$ws.Range("S7").Value = 0.2217647799140164
$ws.Range("T7").Value = 0.2217647799140164
$ws.Range("G8").Value = 1.206946333333333
$ws.Range("H8").Value = 3.620839
$ws.Range("I8").Value = 0.01148164516212593
$ws.Range("J8").Value = 0.01148164516212593
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 0.5174697917923334
$ws.Range("R8").Value = 4.657228126131001
$ws.Range("S8").Value = 0.0000473101642273277269378
$ws.Range("T8").Value = 0.0000473101642273277269378
$ws.Range("G9").Value = 1.206946333333333
$ws.Range("H9").Value = 3.620839
$ws.Range("I9").Value = 0.01148164516212593
$ws.Range("J9").Value = 0.01148164516212593
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 96.82876258428099
$ws.Range("R9").Value = 871.4588632585289
$ws.Range("S9").Value = 0.008852661029592356
$ws.Range("T9").Value = 0.008852661029592358
$ws.Range("G10").Value = 1.206946333333333
$ws.Range("H10").Value = 3.620839
$ws.Range("I10").Value = 0.01148164516212593
$ws.Range("J10").Value = 0.01148164516212593
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 28.23787050148189
$ws.Range("R10").Value = 254.140834513337
$ws.Range("S10").Value = 0.002581673968306248
$ws.Range("T10").Value = 0.002581673968306249
